$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4304569862216283
$ws.Range("C2").Value = 0.1468426154909395
$ws.Range("D2").Value = 0.05358668348728912
$ws.Range("E2").Value = 0.120689926605035
$ws.Range("F2").Value = 1.019680491591998
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 0.9743034593364719
$ws.Range("K2").Value = 0.2759078690219496
$ws.Range("L2").Value = 0.211234619343621
$ws.Range("M2").Value = 0.1391969961823527
$ws.Range("O2").Value = 3.673779994777703
$ws.Range("B3").Value = 0.3949936695926226
$ws.Range("C3").Value = 0.1450407181090654
$ws.Range("D3").Value = 0.05133406264796037
$ws.Range("E3").Value = 0.1212344157329568
$ws.Range("F3").Value = 1.022163523110535
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 0.9820457855585829
$ws.Range("K3").Value = 0.2421641450222012
$ws.Range("L3").Value = 0.2086123896107495
$ws.Range("M3").Value = 0.1323664944757859
$ws.Range("O3").Value = 3.695655788436881
$ws.Range("B4").Value = 0.37330110422144
$ws.Range("C4").Value = 0.1439243604575324
$ws.Range("D4").Value = 0.0499356754888538
$ws.Range("E4").Value = 0.1216141015285643
$ws.Range("F4").Value = 1.024239218371044
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 0.9872413179532558
$ws.Range("K4").Value = 0.2214114467003441
$ws.Range("L4").Value = 0.2071037662165764
$ws.Range("M4").Value = 0.1282237426441846
$ws.Range("O4").Value = 3.710914529281439
$ws.Range("B5").Value = 0.3644824466773002
$ws.Range("C5").Value = 0.1434669499546146
$ws.Range("D5").Value = 0.04936200680949554
$ws.Range("E5").Value = 0.1217802536035411
$ws.Range("F5").Value = 1.025223740669674
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 0.9894696601875523
$ws.Range("K5").Value = 0.2129465139255302
$ws.Range("L5").Value = 0.2065145667823813
$ws.Range("M5").Value = 0.1265485363598735
$ws.Range("O5").Value = 3.717592163080695
$ws.Range("B6").Value = 0.3630194159510154
$ws.Range("C6").Value = 0.1433908479062396
$ws.Range("D6").Value = 0.04926651990708564
$ws.Range("E6").Value = 0.121808533753871
$ws.Range("F6").Value = 1.025395596899742
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 0.9898463875515056
$ws.Range("K6").Value = 0.2115404471210383
$ws.Range("L6").Value = 0.2064182775000063
$ws.Range("M6").Value = 0.1262711583505478
$ws.Range("O6").Value = 3.718728741477392
$ws.Range("B7").Value = 0.373182085921087
$ws.Range("C7").Value = 0.143918201689246
$ws.Range("D7").Value = 0.04992795420477592
$ws.Range("E7").Value = 0.1216162960204965
$ws.Range("F7").Value = 1.024251934449921
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 0.9872709201692018
$ws.Range("K7").Value = 0.2212973175603707
$ws.Range("L7").Value = 0.2070957164145852
$ws.Range("M7").Value = 0.1282010974314929
$ws.Range("O7").Value = 3.71100272515595
$ws.Range("B8").Value = 0.4182125558631924
$ws.Range("C8").Value = 0.1462234083133822
$ws.Range("D8").Value = 0.05281316320583329
$ws.Range("E8").Value = 0.120868261836069
$ws.Range("F8").Value = 1.020422305380549
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 0.9768813649245907
$ws.Range("K8").Value = 0.2642803546539767
$ws.Range("L8").Value = 0.2103094620511072
$ws.Range("M8").Value = 0.1368312832843337
$ws.Range("O8").Value = 3.680943697643826
$ws.Range("B9").Value = 0.5071461540978532
$ws.Range("C9").Value = 0.150663775999206
$ws.Range("D9").Value = 0.05834904085646997
$ws.Range("E9").Value = 0.1197605389680234
$ws.Range("F9").Value = 1.017282644966215
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 0.9600102603098222
$ws.Range("K9").Value = 0.3482834488540334
$ws.Range("L9").Value = 0.2174139395229133
$ws.Range("M9").Value = 0.1541568723307023
$ws.Range("O9").Value = 3.636487554016071
$ws.Range("B10").Value = 0.572846109590472
$ws.Range("C10").Value = 0.1538763405016397
$ws.Range("D10").Value = 0.0623410694214499
$ws.Range("E10").Value = 0.1191646370864952
$ws.Range("F10").Value = 1.017638205870057
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 0.9497479141042788
$ws.Range("K10").Value = 0.4098084395359933
$ws.Range("L10").Value = 0.2231201489644121
$ws.Range("M10").Value = 0.1671263506892799
$ws.Range("O10").Value = 3.612652957088443
$ws.Range("B11").Value = 0.6028085274433579
$ws.Range("C11").Value = 0.1553268360063385
$ws.Range("D11").Value = 0.06414067256979905
$ws.Range("E11").Value = 0.1189406587721411
$ws.Range("F11").Value = 1.01837772825445
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 0.9455419240679177
$ws.Range("K11").Value = 0.4377527216997237
$ws.Range("L11").Value = 0.2258212045567234
$ws.Range("M11").Value = 0.1730778080252264
$ws.Range("O11").Value = 3.603725710979575
$ws.Range("B12").Value = 0.6141648029472719
$ws.Range("C12").Value = 0.1558745098917882
$ws.Range("D12").Value = 0.06481975759362513
$ws.Range("E12").Value = 0.1188625988993568
$ws.Range("F12").Value = 1.018740795757601
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 0.9440156794015095
$ws.Range("K12").Value = 0.4483277895185438
$ws.Range("L12").Value = 0.2268591028728935
$ws.Range("M12").Value = 0.1753387809628748
$ws.Range("O12").Value = 3.600620486997514
$ws.Range("B13").Value = 0.6117185844162805
$ws.Range("C13").Value = 0.1557566299965458
$ws.Range("D13").Value = 0.06467361104114389
$ws.Range("E13").Value = 0.1188791102940385
$ws.Range("F13").Value = 1.018658910971936
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 0.944341427443625
$ws.Range("K13").Value = 0.44605057151594
$ws.Range("L13").Value = 0.2266349039133928
$ws.Range("M13").Value = 0.1748515180672001
$ws.Range("O13").Value = 3.601277008399393
$ws.Range("B14").Value = 0.603742615211786
$ws.Range("C14").Value = 0.1553719256346824
$ws.Range("D14").Value = 0.06419658927937633
$ws.Range("E14").Value = 0.1189341014401162
$ws.Range("F14").Value = 1.018405934393741
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 0.9454150269106165
$ws.Range("K14").Value = 0.4386228787441269
$ws.Range("L14").Value = 0.2259062915129988
$ws.Range("M14").Value = 0.1732636743756899
$ws.Range("O14").Value = 3.603464724752314
$ws.Range("B15").Value = 0.5988584073218419
$ws.Range("C15").Value = 0.1551360742764558
$ws.Range("D15").Value = 0.06390408814574045
$ws.Range("E15").Value = 0.118968664408353
$ws.Range("F15").Value = 1.018261789715154
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 0.9460812936542524
$ws.Range("K15").Value = 0.4340722977313476
$ws.Range("L15").Value = 0.2254619554668551
$ws.Range("M15").Value = 0.1722920191162842
$ws.Range("O15").Value = 3.604840618002925
$ws.Range("B16").Value = 0.5708894520980436
$ws.Range("C16").Value = 0.1537813251886746
$ws.Range("D16").Value = 0.06222312885029879
$ws.Range("E16").Value = 0.1191802209127495
$ws.Range("F16").Value = 1.017601497819555
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 0.9500320884208335
$ws.Range("K16").Value = 0.4079812889665106
$ws.Range("L16").Value = 0.2229457409249846
$ws.Range("M16").Value = 0.1667384355704229
$ws.Range("O16").Value = 3.61327490636981
$ws.Range("B17").Value = 0.5537501943406653
$ws.Range("C17").Value = 0.1529474149506811
$ws.Range("D17").Value = 0.06118769631744669
$ws.Range("E17").Value = 0.119322055335827
$ws.Range("F17").Value = 1.017344372151527
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 0.9525741909230767
$ws.Range("K17").Value = 0.3919637151062432
$ws.Range("L17").Value = 0.2214290372116778
$ws.Range("M17").Value = 0.1633446083924355
$ws.Range("O17").Value = 3.618939551067854
$ws.Range("B18").Value = 0.543899261671271
$ws.Range("C18").Value = 0.1524667459893436
$ws.Range("D18").Value = 0.06059060185557996
$ws.Range("E18").Value = 0.1194080690850079
$ws.Range("F18").Value = 1.017250866750011
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 0.9540798658029672
$ws.Range("K18").Value = 0.3827467447112554
$ws.Range("L18").Value = 0.2205665823250769
$ws.Range("M18").Value = 0.161397431246634
$ws.Range("O18").Value = 3.622377980246455
$ws.Range("B19").Value = 0.5405651470550765
$ws.Range("C19").Value = 0.1523038243892429
$ws.Range("D19").Value = 0.06038817212692038
$ws.Range("E19").Value = 0.1194379539869832
$ws.Range("F19").Value = 1.017228550080816
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 0.954597137489646
$ws.Range("K19").Value = 0.3796253515259025
$ws.Range("L19").Value = 0.2202762746381524
$ws.Range("M19").Value = 0.1607389896648783
$ws.Range("O19").Value = 3.623573138103865
$ws.Range("B20").Value = 0.5555739659356504
$ws.Range("C20").Value = 0.1530362924900643
$ws.Range("D20").Value = 0.06129807957940159
$ws.Range("E20").Value = 0.1193064980060008
$ws.Range("F20").Value = 1.017366115230459
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 0.9522990752703748
$ws.Range("K20").Value = 0.3936692399071546
$ws.Range("L20").Value = 0.2215894674686751
$ws.Range("M20").Value = 0.16370538479309
$ws.Range("O20").Value = 3.618317882699387
$ws.Range("B21").Value = 0.606085080212182
$ws.Range("C21").Value = 0.155484966199964
$ws.Range("D21").Value = 0.06433676710748415
$ws.Range("E21").Value = 0.1189177659920926
$ws.Range("F21").Value = 1.018477986893657
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 0.945097880842134
$ws.Range("K21").Value = 0.4408047616568354
$ws.Range("L21").Value = 0.2261198942833715
$ws.Range("M21").Value = 0.1737298658864717
$ws.Range("O21").Value = 3.602814667666422
$ws.Range("B22").Value = 0.6391559208910849
$ws.Range("C22").Value = 0.1570759931766474
$ws.Range("D22").Value = 0.06630880648133086
$ws.Range("E22").Value = 0.1187030774922526
$ws.Range("F22").Value = 1.019688573622261
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 0.940778926942258
$ws.Range("K22").Value = 0.4715705190785684
$ws.Range("L22").Value = 0.2291685801826873
$ws.Range("M22").Value = 0.180323846607358
$ws.Range("O22").Value = 3.594287219584032
$ws.Range("B23").Value = 0.6215002160051881
$ws.Range("C23").Value = 0.156227694234687
$ws.Range("D23").Value = 0.06525757497637841
$ws.Range("E23").Value = 0.1188140640541331
$ws.Range("F23").Value = 1.018998201635249
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 0.9430485896413145
$ws.Range("K23").Value = 0.4551540945910233
$ws.Range("L23").Value = 0.2275334297957272
$ws.Range("M23").Value = 0.1768006785427616
$ws.Range("O23").Value = 3.598691665875464
$ws.Range("B24").Value = 0.5547494304790348
$ws.Range("C24").Value = 0.1529961148249015
$ws.Range("D24").Value = 0.06124818094369999
$ws.Range("E24").Value = 0.1193135175462103
$ws.Range("F24").Value = 1.017356115978465
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 0.9524233174174128
$ws.Range("K24").Value = 0.3928981978151
$ws.Range("L24").Value = 0.2215169072961345
$ws.Range("M24").Value = 0.1635422653777141
$ws.Range("O24").Value = 3.618598372822731
$ws.Range("B25").Value = 0.4830223065543748
$ws.Range("C25").Value = 0.14947121379344
$ws.Range("D25").Value = 0.05686458404044714
$ws.Range("E25").Value = 0.1200218609907555
$ws.Range("F25").Value = 1.017664340423295
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 0.964199606627048
$ws.Range("K25").Value = 0.3255908650441199
$ws.Range("L25").Value = 0.2154063605007508
$ws.Range("M25").Value = 0.1494272990656142
$ws.Range("O25").Value = 3.646963519823373
